# Daily attendance processing - rotate "Recorded By" (column G) value lists
# so that the most recently appended recorder (last in the comma-separated
# list) is moved to the front of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the used range extents
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Column G is the 7th column ("Recorded By"); data starts at row 2 (row 1 is header)
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -eq "") { continue }

    $parts = $val -split ", "

    if ($parts.Count -gt 1) {
        $rotated = @($parts[$parts.Count - 1]) + $parts[0..($parts.Count - 2)]
        $newVal = [string]::Join(", ", $rotated)
        $cell.Value = $newVal
    }
}
